# Append the latest batch of combined-store data onto Sheet1's table.
# The sheet holds a small two-column table (header row + 5 data rows);
# this refreshes it with the newest aggregated counts while keeping the
# same "combinedstore1" / "combinedstore2" column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Header row (unchanged labels, refreshed as part of the same write).
$ws.Range("A1").Value = "combinedstore1"
$ws.Range("B1").Value = "combinedstore2"

# Newest rows appended/merged onto the table.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 9

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 8

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 7

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 7

$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 6
